$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(18).Insert()
$ws.Cells.Item(18, 1).Value = 10
$ws.Cells.Item(18, 2).Value = 2023
$ws.Cells.Item(18, 3).Value = "Reunión del grupo Almería"

$ws.Cells.Item(41, 1).Value = 6
$ws.Cells.Item(41, 2).Value = 2025
$ws.Cells.Item(41, 3).Value = "Reunión del grupo Pontevedra"
